$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing rows 66-123 down to 67-124
$ws.Rows("66:66").Insert()

# Populate the newly inserted row 66 with the new price-report entry
$ws.Range("A66").Value = 11
$ws.Range("B66").Value = "Vega Monumental Concepción"
$ws.Range("C66").Value = "Bíobío"
$ws.Range("D66").Value = 44778
$ws.Range("E66").Value = 8
$ws.Range("F66").Value = 100112021
$ws.Range("G66").Value = "Ají"
$ws.Range("H66").Value = "Inferno"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 22
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 16000
$ws.Range("M66").Value = 15455
$ws.Range("N66").Value = "`$/caja 12 kilos"
$ws.Range("O66").Value = "Región de Arica y Parinacota"
$ws.Range("P66").Value = 1288
$ws.Range("Q66").Value = 12
$ws.Range("R66").Value = "Hortaliza"
